$wb = $excel.ActiveWorkbook

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 311
$ws.Range("I41").Value = 360.625
$ws.Range("J41").Value = 266.8889
$ws.Range("K41").Value = 360.625
$ws.Range("L41").Value = 266.8889
$ws.Range("M41").Value = 79.375
$ws.Range("N41").Value = -1146.8889

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3382.75
$ws.Range("I64").Value = 3194.4666
$ws.Range("J64").Value = 3600
$ws.Range("K64").Value = 3194.4666
$ws.Range("L64").Value = 3600
$ws.Range("M64").Value = -2946.4666
$ws.Range("N64").Value = -4096

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3382.75
$ws.Range("I67").Value = 3194.4666
$ws.Range("J67").Value = 3600
$ws.Range("K67").Value = 3194.4666
$ws.Range("L67").Value = 3600
$ws.Range("M67").Value = -2336.4666
$ws.Range("N67").Value = -5316

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 685.1579
$ws.Range("I135").Value = 511.31033
$ws.Range("J135").Value = 1245.3334
$ws.Range("K135").Value = 4601.79297
$ws.Range("L135").Value = 11208.0006
$ws.Range("M135").Value = -2066.79297
$ws.Range("N135").Value = -16278.0006

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1668.9166
$ws.Range("I137").Value = 1425.4783
$ws.Range("J137").Value = 2099.6155
$ws.Range("K137").Value = 4276.4349
$ws.Range("L137").Value = 6298.8465
$ws.Range("M137").Value = -1726.4349
$ws.Range("N137").Value = -11398.8465

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2535339.8
$ws.Range("I138").Value = 14287732
$ws.Range("J138").Value = 4055.2615
$ws.Range("K138").Value = 42863196
$ws.Range("L138").Value = 12165.7845
$ws.Range("M138").Value = -42858056
$ws.Range("N138").Value = -22445.7845

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1465.7037
$ws.Range("I2").Value = 960.06665
$ws.Range("J2").Value = 2097.75
$ws.Range("K2").Value = 960.06665
$ws.Range("L2").Value = 2097.75
$ws.Range("M2").Value = -847.06665
$ws.Range("N2").Value = -2323.75

# ARM row 16
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 3503
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 6
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 281
$ws.Range("N16").Value = -7574

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16251
$ws.Range("I32").Value = 16984.762
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 16984.762
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -16697.762
$ws.Range("N32").Value = -5574

# ARM row 58
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 45762.145
$ws.Range("J58").Value = 45762.145
$ws.Range("L58").Value = 45762.145
$ws.Range("N58").Value = -46622.145

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2552
$ws.Range("I61").Value = 2062.4
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2062.4
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1850.4
$ws.Range("N61").Value = -5424

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3435.318
$ws.Range("I63").Value = 3829
$ws.Range("J63").Value = 2866.6667
$ws.Range("K63").Value = 3829
$ws.Range("L63").Value = 2866.6667
$ws.Range("M63").Value = -3143
$ws.Range("N63").Value = -4238.6667

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3435.318
$ws.Range("I66").Value = 3829
$ws.Range("J66").Value = 2866.6667
$ws.Range("K66").Value = 19145
$ws.Range("L66").Value = 14333.3335
$ws.Range("M66").Value = -15713
$ws.Range("N66").Value = -21197.3335

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1465.7037
$ws.Range("I116").Value = 960.06665
$ws.Range("J116").Value = 2097.75
$ws.Range("K116").Value = 960.06665
$ws.Range("L116").Value = 2097.75
$ws.Range("M116").Value = 1333.93335
$ws.Range("N116").Value = -6685.75

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2179.0454
$ws.Range("I132").Value = 1619.4546
$ws.Range("J132").Value = 2738.6365
$ws.Range("K132").Value = 4858.3638
$ws.Range("L132").Value = 8215.9095
$ws.Range("M132").Value = -2328.3638
$ws.Range("N132").Value = -13275.9095

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2552
$ws.Range("I136").Value = 2062.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6187.200000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3637.200000000001
$ws.Range("N136").Value = -20100

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1465.7037
$ws.Range("I3").Value = 960.06665
$ws.Range("J3").Value = 2097.75
$ws.Range("K3").Value = 960.06665
$ws.Range("L3").Value = 2097.75
$ws.Range("M3").Value = -846.06665
$ws.Range("N3").Value = -2325.75

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 38226.145
$ws.Range("I107").Value = 52876.6
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 52876.6
$ws.Range("L107").Value = 1600
$ws.Range("M107").Value = -50956.6
$ws.Range("N107").Value = -5440

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2730.5557
$ws.Range("I134").Value = 2211.5386
$ws.Range("J134").Value = 4080
$ws.Range("K134").Value = 6634.6158
$ws.Range("L134").Value = 12240
$ws.Range("M134").Value = -4099.6158
$ws.Range("N134").Value = -17310

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3656.88
$ws.Range("I99").Value = 3715.2104
$ws.Range("K99").Value = 3715.2104
$ws.Range("M99").Value = -2217.2104

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3656.88
$ws.Range("I126").Value = 3715.2104
$ws.Range("K126").Value = 11145.6312
$ws.Range("M126").Value = -8675.6312

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2375
$ws.Range("I132").Value = 1791.8889
$ws.Range("K132").Value = 5375.6667
$ws.Range("M132").Value = -2845.6667

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1431.85
$ws.Range("I134").Value = 1358.5625
$ws.Range("J134").Value = 1725
$ws.Range("K134").Value = 4075.6875
$ws.Range("L134").Value = 5175
$ws.Range("M134").Value = -1540.6875
$ws.Range("N134").Value = -10245

# CRP row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

# CUL row 31
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1933.6666
$ws.Range("I31").Value = 400.5
$ws.Range("K31").Value = 1201.5
$ws.Range("M31").Value = -913.5

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 15225
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 15225
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2442.361
$ws.Range("I140").Value = 953.2778
$ws.Range("K140").Value = 2859.8334
$ws.Range("M140").Value = 2320.1666

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2147639.8
$ws.Range("I14").Value = 7868833.5
$ws.Range("J14").Value = 2192.1875
$ws.Range("K14").Value = 7868833.5
$ws.Range("L14").Value = 2192.1875
$ws.Range("M14").Value = -7868665.5
$ws.Range("N14").Value = -2528.1875

# GSM row 74
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 29333.334
$ws.Range("J74").Value = 29333.334
$ws.Range("L74").Value = 29333.334
$ws.Range("N74").Value = -31205.334

# GSM row 77
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 29333.334
$ws.Range("J77").Value = 29333.334
$ws.Range("L77").Value = 88000.00199999999
$ws.Range("N77").Value = -97360.00199999999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2729.0557
$ws.Range("I132").Value = 2093.8333
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 6281.499899999999
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -3751.499899999999
$ws.Range("N132").Value = -17058.5

# WVR row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2459.6667
$ws.Range("I23").Value = 2459.6667
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2459.6667
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -2230.6667
$ws.Range("N23").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1403.6428
$ws.Range("I132").Value = 1101
$ws.Range("J132").Value = 2042.5555
$ws.Range("K132").Value = 3303
$ws.Range("L132").Value = 6127.666499999999
$ws.Range("M132").Value = -773
$ws.Range("N132").Value = -11187.6665

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1248.5238
$ws.Range("I136").Value = 1271.7059
$ws.Range("J136").Value = 1150
$ws.Range("K136").Value = 3815.1177
$ws.Range("L136").Value = 3450
$ws.Range("M136").Value = -1265.1177
$ws.Range("N136").Value = -8550
